$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3318.0908
$ws.Range("I40").Value = 6333
$ws.Range("J40").Value = 2187.5
$ws.Range("K40").Value = 6333
$ws.Range("L40").Value = 2187.5
$ws.Range("M40").Value = -6158
$ws.Range("N40").Value = -2537.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1323.3334
$ws.Range("I112").Value = 5550
$ws.Range("J112").Value = 1169.6364
$ws.Range("K112").Value = 16650
$ws.Range("L112").Value = 3508.9092
$ws.Range("M112").Value = -15542
$ws.Range("N112").Value = -5724.9092

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6253386.5
$ws.Range("I137").Value = 2382.1875
$ws.Range("J137").Value = 12504390
$ws.Range("K137").Value = 7146.5625
$ws.Range("L137").Value = 37513170
$ws.Range("M137").Value = -4596.5625
$ws.Range("N137").Value = -37518270

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8069406.5
$ws.Range("I138").Value = 2327
$ws.Range("J138").Value = 27788934
$ws.Range("K138").Value = 6981
$ws.Range("L138").Value = 83366802
$ws.Range("M138").Value = -1841
$ws.Range("N138").Value = -83377082

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2768.842
$ws.Range("I45").Value = 2831.3845
$ws.Range("J45").Value = 2633.3333
$ws.Range("K45").Value = 2831.3845
$ws.Range("L45").Value = 2633.3333
$ws.Range("M45").Value = -2454.3845
$ws.Range("N45").Value = -3387.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4809411
$ws.Range("I132").Value = 6580182
$ws.Range("J132").Value = 3032.0715
$ws.Range("K132").Value = 19740546
$ws.Range("L132").Value = 9096.2145
$ws.Range("M132").Value = -19738016
$ws.Range("N132").Value = -14156.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 39066.125
$ws.Range("J135").Value = 39066.125
$ws.Range("L135").Value = 39066.125
$ws.Range("N135").Value = -49206.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 46880.832
$ws.Range("J139").Value = 46880.832
$ws.Range("L139").Value = 46880.832
$ws.Range("N139").Value = -57160.832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 56018.5
$ws.Range("J55").Value = 56018.5
$ws.Range("L55").Value = 56018.5
$ws.Range("N55").Value = -56564.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13890553
$ws.Range("I86").Value = 1621.1538
$ws.Range("J86").Value = 50001776
$ws.Range("K86").Value = 1621.1538
$ws.Range("L86").Value = 50001776
$ws.Range("M86").Value = -498.1538
$ws.Range("N86").Value = -50004022

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 13890553
$ws.Range("I89").Value = 1621.1538
$ws.Range("J89").Value = 50001776
$ws.Range("K89").Value = 8105.769
$ws.Range("L89").Value = 250008880
$ws.Range("M89").Value = -2489.769
$ws.Range("N89").Value = -250020112

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2822.9307
$ws.Range("I134").Value = 1967.2181
$ws.Range("J134").Value = 5591.4116
$ws.Range("K134").Value = 5901.6543
$ws.Range("L134").Value = 16774.2348
$ws.Range("M134").Value = -3366.6543
$ws.Range("N134").Value = -21844.2348

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 43916.332
$ws.Range("J138").Value = 43916.332
$ws.Range("L138").Value = 43916.332
$ws.Range("N138").Value = -54196.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6292685.5
$ws.Range("I31").Value = 5166.0347
$ws.Range("J31").Value = 13890105
$ws.Range("K31").Value = 5166.0347
$ws.Range("L31").Value = 13890105
$ws.Range("M31").Value = -4871.0347
$ws.Range("N31").Value = -13890695

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6292685.5
$ws.Range("I34").Value = 5166.0347
$ws.Range("J34").Value = 13890105
$ws.Range("K34").Value = 5166.0347
$ws.Range("L34").Value = 13890105
$ws.Range("M34").Value = -4964.0347
$ws.Range("N34").Value = -13890509

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3772.611
$ws.Range("I58").Value = 2516
$ws.Range("J58").Value = 4572.273
$ws.Range("K58").Value = 2516
$ws.Range("L58").Value = 4572.273
$ws.Range("M58").Value = -2313
$ws.Range("N58").Value = -4978.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2365.5557
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 2496.6667
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 2496.6667
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -3744.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2365.5557
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 2496.6667
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 12483.3335
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -18723.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3772.611
$ws.Range("I136").Value = 2516
$ws.Range("J136").Value = 4572.273
$ws.Range("K136").Value = 7548
$ws.Range("L136").Value = 13716.819
$ws.Range("M136").Value = -4998
$ws.Range("N136").Value = -18816.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5555709.5
$ws.Range("I4").Value = 6250162
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 18750486
$ws.Range("L4").Value = 270
$ws.Range("M4").Value = -18750374
$ws.Range("N4").Value = -494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 835.2292
$ws.Range("I68").Value = 548.9079
$ws.Range("J68").Value = 1923.25
$ws.Range("K68").Value = 1646.7237
$ws.Range("L68").Value = 5769.75
$ws.Range("M68").Value = -835.7237
$ws.Range("N68").Value = -7391.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 835.2292
$ws.Range("I71").Value = 548.9079
$ws.Range("J71").Value = 1923.25
$ws.Range("K71").Value = 4940.1711
$ws.Range("L71").Value = 17309.25
$ws.Range("M71").Value = -884.1711000000005
$ws.Range("N71").Value = -25421.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1234.125
$ws.Range("I121").Value = 129.8
$ws.Range("K121").Value = 389.4
$ws.Range("M121").Value = 920.5999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2339.0938
$ws.Range("I102").Value = 2525.8076
$ws.Range("K102").Value = 2525.8076
$ws.Range("M102").Value = -903.8076000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3289.3333
$ws.Range("I132").Value = 1813.409
$ws.Range("J132").Value = 5199.353
$ws.Range("K132").Value = 5440.227000000001
$ws.Range("L132").Value = 15598.059
$ws.Range("M132").Value = -2910.227000000001
$ws.Range("N132").Value = -20658.059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5829.615
$ws.Range("I122").Value = 5624.4116
$ws.Range("K122").Value = 16873.2348
$ws.Range("M122").Value = -14423.2348

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7468200.5
$ws.Range("I132").Value = 3672.1025
$ws.Range("J132").Value = 17865222
$ws.Range("K132").Value = 11016.3075
$ws.Range("L132").Value = 53595666
$ws.Range("M132").Value = -8486.307499999999
$ws.Range("N132").Value = -53600726

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 34888428
$ws.Range("I136").Value = 62501708
$ws.Range("K136").Value = 187505124
$ws.Range("M136").Value = -187502574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7581.4287
$ws.Range("I74").Value = 4569
$ws.Range("J74").Value = 8083.5
$ws.Range("K74").Value = 4569
$ws.Range("L74").Value = 8083.5
$ws.Range("M74").Value = -3633
$ws.Range("N74").Value = -9955.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 7581.4287
$ws.Range("I77").Value = 4569
$ws.Range("J77").Value = 8083.5
$ws.Range("K77").Value = 13707
$ws.Range("L77").Value = 24250.5
$ws.Range("M77").Value = -9027
$ws.Range("N77").Value = -33610.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3282.4695
$ws.Range("I132").Value = 3759.838
$ws.Range("J132").Value = 1810.5834
$ws.Range("K132").Value = 11279.514
$ws.Range("L132").Value = 5431.7502
$ws.Range("M132").Value = -8749.514000000001
$ws.Range("N132").Value = -10491.7502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 39273
$ws.Range("J133").Value = 39273
$ws.Range("L133").Value = 39273
$ws.Range("N133").Value = -49393
